$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after existing row 226, pushing the
# remaining data (old rows 227-279) down to rows 229-281.
$ws.Rows("227:228").Insert()

# New row 227 (fresh weekly record, not shifted from elsewhere).
$ws.Range("A227").Value = 7
$ws.Range("B227").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C227").Value = "Ñuble"
$ws.Range("D227").Value = 44754
$ws.Range("E227").Value = 16
$ws.Range("F227").Value = 100112008
$ws.Range("G227").Value = "Coliflor"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 200
$ws.Range("K227").Value = 1000
$ws.Range("L227").Value = 1200
$ws.Range("M227").Value = 1100
$ws.Range("N227").Value = "$/unidad"
$ws.Range("O227").Value = "Provincia de Diguillín"
$ws.Range("P227").Value = 1100
$ws.Range("Q227").Value = 1
$ws.Range("R227").Value = "Hortaliza"

# New row 228 (fresh weekly record, not shifted from elsewhere).
$ws.Range("A228").Value = 7
$ws.Range("B228").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C228").Value = "Ñuble"
$ws.Range("D228").Value = 44754
$ws.Range("E228").Value = 16
$ws.Range("F228").Value = 100112008
$ws.Range("G228").Value = "Coliflor"
$ws.Range("H228").Value = "Sin especificar"
$ws.Range("I228").Value = "Segunda"
$ws.Range("J228").Value = 120
$ws.Range("K228").Value = 900
$ws.Range("L228").Value = 900
$ws.Range("M228").Value = 900
$ws.Range("N228").Value = "$/unidad"
$ws.Range("O228").Value = "Provincia de Diguillín"
$ws.Range("P228").Value = 900
$ws.Range("Q228").Value = 1
$ws.Range("R228").Value = "Hortaliza"
